$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The dataset gained two new gyroscope samples at the top (the "falling"
# class rows being folded into this sheet). That shifts every existing
# reading down by two rows, and the two oldest readings (previously rows
# 21 and 22) fall off the end, so the used range shrinks from A1:C22 to
# A1:C21 overall (net: +2 new rows, -3 old rows).
#
# Rebuild rows 2-21 in place with the new, shifted x/y/z values.
$data = @(
  @(0.0204639863222837, -0.0009162978967650999, 0.007177666760981),
  @(0.0155770638957619, -0.046578474342823, 0.0195476878434419),
  @(0.0387899428606033, -0.2247984111309051, 0.1269072592258453),
  @(-0.1467603743076324, -0.3396410644054413, 0.0708603709936142),
  @(0.0829249545931816, -0.4129449129104614, 0.1215621829032898),
  @(0.2186897546052932, -0.2243402600288391, 0.5036584138870239),
  @(0.2626720666885376, 0.047036625444889, 0.1678352206945419),
  @(-0.2110539376735687, 0.9388998746871948, -0.7550294399261475),
  @(0.06963863968849179, 1.037860035896301, -2.770426511764526),
  @(-0.3165809214115143, -1.475086808204651, -1.813048005104065),
  @(-2.540435791015625, -5.590180397033691, 0.3164281845092773),
  @(0.2667953968048095, 0.7530441284179688, 0.5949827432632446),
  @(0.5946773290634155, 0.5623014569282532, 0.265573650598526),
  @(-0.2397646158933639, -0.2414444833993911, -0.8886561989784241),
  @(0.3083342313766479, 0.1081231459975242, 0.2041817009449005),
  @(0.0649044290184974, 0.2634356319904327, 0.0945313945412635),
  @(-0.0543670058250427, 0.1363756656646728, 0.0697913542389869),
  @(-0.08918632566928859, -0.3081815242767334, -0.1731802970170974),
  @(-0.0502436682581901, -0.1233947798609733, -0.0299323964864015),
  @(0.0731511116027832, 0.2703078687191009, 0.0377209298312664)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

# Drop the now-unused trailing row so the sheet ends at row 21 (A1:C21).
$ws.Range("A22:C22").EntireRow.Delete()
